# Refresh cryptos list figures (price + 1h volume change) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Ws, [string]$CellRef, [string]$Text, [bool]$ForceText = $false)
    # When $ForceText is set, prefix with an apostrophe so Excel keeps the
    # numeric-looking string as literal text (preserves trailing zeros, etc.)
    if ($ForceText) {
        $Ws.Range($CellRef).Value = "'" + $Text
    } else {
        $Ws.Range($CellRef).Value = $Text
    }
}

Set-TextValue $ws 'D2' '30.466.34' $false
Set-TextValue $ws 'E2' '  +0.91%  ' $false
Set-TextValue $ws 'D3' '1.878.93' $false
Set-TextValue $ws 'E4' '  -0.02%  ' $false
Set-TextValue $ws 'D5' '246.73' $true
Set-TextValue $ws 'E5' '  +5.45%  ' $false
Set-TextValue $ws 'E7' '  +2.18%  ' $false
Set-TextValue $ws 'D8' '0.2902' $true
Set-TextValue $ws 'E8' '  +1.96%  ' $false
Set-TextValue $ws 'D9' '0.06523' $true
Set-TextValue $ws 'E9' '  +0.55%  ' $false
Set-TextValue $ws 'D10' '21.61' $true
Set-TextValue $ws 'E10' '  +1.43%  ' $false
Set-TextValue $ws 'D11' '0.07741' $true
Set-TextValue $ws 'E11' '  +0.03%  ' $false
Set-TextValue $ws 'D12' '0.7447' $true
Set-TextValue $ws 'E12' '  +8.98%  ' $false
Set-TextValue $ws 'E13' '  +3.40%  ' $false
Set-TextValue $ws 'D14' '1.873.93' $false
Set-TextValue $ws 'E14' '  +0.69%  ' $false
Set-TextValue $ws 'D15' '5.120' $true
Set-TextValue $ws 'E15' '  +1.53%  ' $false
Set-TextValue $ws 'D16' '273.84' $true
Set-TextValue $ws 'E16' '  +1.09%  ' $false
Set-TextValue $ws 'D17' '30.459.38' $false
Set-TextValue $ws 'E18' '  +2.29%  ' $false
Set-TextValue $ws 'D19' '0.000007573' $true
Set-TextValue $ws 'E19' '  -0.33%  ' $false
Set-TextValue $ws 'E20' '  -0.01%  ' $false
Set-TextValue $ws 'D21' '2.125.71' $false
Set-TextValue $ws 'E21' '  +1.26%  ' $false
Set-TextValue $ws 'E22' '  +0.03%  ' $false
Set-TextValue $ws 'D23' '5.247' $true
Set-TextValue $ws 'E23' '  +2.10%  ' $false
Set-TextValue $ws 'D24' '6.173' $true
Set-TextValue $ws 'E24' '  +1.30%  ' $false
Set-TextValue $ws 'D25' '9.288' $true
Set-TextValue $ws 'E25' '  -0.48%  ' $false
Set-TextValue $ws 'D26' '164.56' $true
Set-TextValue $ws 'E26' '  -0.57%  ' $false
Set-TextValue $ws 'D27' '18.90' $true
Set-TextValue $ws 'D28' '1.957' $true
Set-TextValue $ws 'E28' '  +3.52%  ' $false
Set-TextValue $ws 'E29' '  +0.60%  ' $false
Set-TextValue $ws 'D30' '0.09997' $true
Set-TextValue $ws 'E30' '  +1.70%  ' $false
Set-TextValue $ws 'D31' '1.515' $true
Set-TextValue $ws 'E31' '  +4.55%  ' $false
Set-TextValue $ws 'D32' '4.327' $true
Set-TextValue $ws 'E32' '  +2.29%  ' $false
Set-TextValue $ws 'D33' '4.065' $true
Set-TextValue $ws 'D34' '0.04773' $true
Set-TextValue $ws 'E34' '  +2.28%  ' $false
Set-TextValue $ws 'E35' '  +0.34%  ' $false
Set-TextValue $ws 'D36' '0.6989' $true
Set-TextValue $ws 'E36' '  +1.27%  ' $false
Set-TextValue $ws 'E37' '  +0.23%  ' $false
Set-TextValue $ws 'D38' '0.01865' $true
Set-TextValue $ws 'E38' '  +1.51%  ' $false
Set-TextValue $ws 'D39' '2.737' $true
Set-TextValue $ws 'E39' '  -0.52%  ' $false
Set-TextValue $ws 'D40' '6.362' $true
Set-TextValue $ws 'E40' '  +0.92%  ' $false
Set-TextValue $ws 'E41' '  +2.72%  ' $false
Set-TextValue $ws 'D42' '70.03' $true
Set-TextValue $ws 'E42' '  -1.03%  ' $false
Set-TextValue $ws 'D43' '0.4170' $true
Set-TextValue $ws 'E43' '  +2.75%  ' $false
Set-TextValue $ws 'D44' '0.9999' $true
Set-TextValue $ws 'E44' '  -0.02%  ' $false
Set-TextValue $ws 'D45' '0.8382' $true
Set-TextValue $ws 'E45' '  +0.64%  ' $false
Set-TextValue $ws 'D46' '102.72' $true
Set-TextValue $ws 'E46' '  +0.56%  ' $false
Set-TextValue $ws 'D47' '9.285' $true
Set-TextValue $ws 'E47' '  +2.37%  ' $false
Set-TextValue $ws 'D48' '7.079' $true
Set-TextValue $ws 'E48' '  +1.83%  ' $false
Set-TextValue $ws 'D49' '35.34' $true
Set-TextValue $ws 'E49' '  +3.87%  ' $false
Set-TextValue $ws 'D50' '923.74' $true
Set-TextValue $ws 'E50' '  -1.15%  ' $false
Set-TextValue $ws 'E51' '  +0.65%  ' $false
